$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Subject" column values for the two log rows are being corrected from
# "Internal Medicine 1" to the correct rotation name "Immuno&hema".
$rng = $ws.Range("B2:B3")
$rng.Value = "Immuno&hema"

# Highlight the corrected column with a light gray fill, centered text and
# the standard 11pt Calibri font.
$rng.Font.Size = 11
$rng.Interior.Color = 15790320     # RGB(240,240,240) -> FFF0F0F0
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.VerticalAlignment = -4108     # xlCenter
